# delete samples for running faster
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 10 and 9 first (bottom-up) so row indices above are unaffected,
# then rows 3 and 2 (also bottom-up) to remove the remaining two sample rows.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(2).Delete()
